# Auto-generated edit script: refreshes scheduled market-data columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several sheets,
# matching the upstream scheduled-runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(16, 8).Value = 0  # H16: 2824.75 -> 0
$ws.Cells.Item(16, 9).Value = 0  # I16: 2824.75 -> 0
$ws.Cells.Item(16, 11).Value = 0  # K16: 2824.75 -> 0
$ws.Cells.Item(16, 13).ClearContents()  # M16: -2594.75 -> (removed)
$ws.Cells.Item(41, 8).Value = 455.57574  # H41: 735.1429000000001 -> 455.57574
$ws.Cells.Item(41, 9).Value = 336.08334  # I41: 567 -> 336.08334
$ws.Cells.Item(41, 10).Value = 523.8570999999999  # J41: 781 -> 523.8570999999999
$ws.Cells.Item(41, 11).Value = 336.08334  # K41: 567 -> 336.08334
$ws.Cells.Item(41, 12).Value = 523.8570999999999  # L41: 781 -> 523.8570999999999
$ws.Cells.Item(41, 13).Value = 103.91666  # M41: -127 -> 103.91666
$ws.Cells.Item(41, 14).Value = -1403.8571  # N41: -1661 -> -1403.8571
$ws.Cells.Item(74, 8).Value = 3075.2122  # H74: 3449.95 -> 3075.2122
$ws.Cells.Item(74, 9).Value = 2964.2727  # I74: 3773 -> 2964.2727
$ws.Cells.Item(74, 10).Value = 3297.0908  # J74: 3276 -> 3297.0908
$ws.Cells.Item(74, 11).Value = 2964.2727  # K74: 3773 -> 2964.2727
$ws.Cells.Item(74, 12).Value = 3297.0908  # L74: 3276 -> 3297.0908
$ws.Cells.Item(74, 13).Value = -2028.2727  # M74: -2837 -> -2028.2727
$ws.Cells.Item(74, 14).Value = -5169.0908  # N74: -5148 -> -5169.0908
$ws.Cells.Item(77, 8).Value = 3075.2122  # H77: 3449.95 -> 3075.2122
$ws.Cells.Item(77, 9).Value = 2964.2727  # I77: 3773 -> 2964.2727
$ws.Cells.Item(77, 10).Value = 3297.0908  # J77: 3276 -> 3297.0908
$ws.Cells.Item(77, 11).Value = 14821.3635  # K77: 18865 -> 14821.3635
$ws.Cells.Item(77, 12).Value = 16485.454  # L77: 16380 -> 16485.454
$ws.Cells.Item(77, 13).Value = -10141.3635  # M77: -14185 -> -10141.3635
$ws.Cells.Item(77, 14).Value = -25845.454  # N77: -25740 -> -25845.454
$ws.Cells.Item(86, 8).Value = 3103.4062  # H86: 3050.2646 -> 3103.4062
$ws.Cells.Item(86, 9).Value = 6100.25  # I86: 4550.1665 -> 6100.25
$ws.Cells.Item(86, 10).Value = 2675.2856  # J86: 2728.8572 -> 2675.2856
$ws.Cells.Item(86, 11).Value = 6100.25  # K86: 4550.1665 -> 6100.25
$ws.Cells.Item(86, 12).Value = 2675.2856  # L86: 2728.8572 -> 2675.2856
$ws.Cells.Item(86, 13).Value = -4977.25  # M86: -3427.1665 -> -4977.25
$ws.Cells.Item(86, 14).Value = -4921.2856  # N86: -4974.8572 -> -4921.2856
$ws.Cells.Item(89, 8).Value = 3103.4062  # H89: 3050.2646 -> 3103.4062
$ws.Cells.Item(89, 9).Value = 6100.25  # I89: 4550.1665 -> 6100.25
$ws.Cells.Item(89, 10).Value = 2675.2856  # J89: 2728.8572 -> 2675.2856
$ws.Cells.Item(89, 11).Value = 30501.25  # K89: 22750.8325 -> 30501.25
$ws.Cells.Item(89, 12).Value = 13376.428  # L89: 13644.286 -> 13376.428
$ws.Cells.Item(89, 13).Value = -24885.25  # M89: -17134.8325 -> -24885.25
$ws.Cells.Item(89, 14).Value = -24608.428  # N89: -24876.286 -> -24608.428
$ws.Cells.Item(129, 8).Value = 598.875  # H129: 603.9375 -> 598.875
$ws.Cells.Item(129, 10).Value = 912.5  # J129: 922.625 -> 912.5
$ws.Cells.Item(129, 12).Value = 2737.5  # L129: 2767.875 -> 2737.5
$ws.Cells.Item(129, 14).Value = -12737.5  # N129: -12767.875 -> -12737.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 24531.857  # H32: 17423.46 -> 24531.857
$ws.Cells.Item(32, 9).Value = 5927.1284  # I32: 5183.284 -> 5927.1284
$ws.Cells.Item(32, 10).Value = 86547.62  # J32: 69605.266 -> 86547.62
$ws.Cells.Item(32, 11).Value = 5927.1284  # K32: 5183.284 -> 5927.1284
$ws.Cells.Item(32, 12).Value = 86547.62  # L32: 69605.266 -> 86547.62
$ws.Cells.Item(32, 13).Value = -5640.1284  # M32: -4896.284 -> -5640.1284
$ws.Cells.Item(32, 14).Value = -87121.62  # N32: -70179.266 -> -87121.62
$ws.Cells.Item(45, 8).Value = 1779.8  # H45: 1778.4286 -> 1779.8
$ws.Cells.Item(45, 9).Value = 1672  # I45: 1616.5 -> 1672
$ws.Cells.Item(45, 11).Value = 1672  # K45: 1616.5 -> 1672
$ws.Cells.Item(45, 13).Value = -1295  # M45: -1239.5 -> -1295
$ws.Cells.Item(61, 8).Value = 3134.5557  # H61: 3585.5715 -> 3134.5557
$ws.Cells.Item(61, 9).Value = 2078  # I61: 1999.5 -> 2078
$ws.Cells.Item(61, 10).Value = 3979.8  # J61: 4220 -> 3979.8
$ws.Cells.Item(61, 11).Value = 2078  # K61: 1999.5 -> 2078
$ws.Cells.Item(61, 12).Value = 3979.8  # L61: 4220 -> 3979.8
$ws.Cells.Item(61, 13).Value = -1866  # M61: -1787.5 -> -1866
$ws.Cells.Item(61, 14).Value = -4403.8  # N61: -4644 -> -4403.8
$ws.Cells.Item(74, 8).Value = 1274.2273  # H74: 1234.52 -> 1274.2273
$ws.Cells.Item(74, 9).Value = 780.82355  # I74: 742.3158 -> 780.82355
$ws.Cells.Item(74, 10).Value = 2951.8  # J74: 2793.1667 -> 2951.8
$ws.Cells.Item(74, 11).Value = 780.82355  # K74: 742.3158 -> 780.82355
$ws.Cells.Item(74, 12).Value = 2951.8  # L74: 2793.1667 -> 2951.8
$ws.Cells.Item(74, 13).Value = 93.17645000000005  # M74: 131.6842 -> 93.17645000000005
$ws.Cells.Item(74, 14).Value = -4699.8  # N74: -4541.1667 -> -4699.8
$ws.Cells.Item(77, 8).Value = 1274.2273  # H77: 1234.52 -> 1274.2273
$ws.Cells.Item(77, 9).Value = 780.82355  # I77: 742.3158 -> 780.82355
$ws.Cells.Item(77, 10).Value = 2951.8  # J77: 2793.1667 -> 2951.8
$ws.Cells.Item(77, 11).Value = 3904.11775  # K77: 3711.579 -> 3904.11775
$ws.Cells.Item(77, 12).Value = 14759  # L77: 13965.8335 -> 14759
$ws.Cells.Item(77, 13).Value = 463.8822500000001  # M77: 656.4210000000003 -> 463.8822500000001
$ws.Cells.Item(77, 14).Value = -23495  # N77: -22701.8335 -> -23495
$ws.Cells.Item(102, 8).Value = 1258.3334  # H102: 799.375 -> 1258.3334
$ws.Cells.Item(102, 9).Value = 1258.3334  # I102: 770.7143 -> 1258.3334
$ws.Cells.Item(102, 10).Value = 0  # J102: 1000 -> 0
$ws.Cells.Item(102, 11).Value = 1258.3334  # K102: 770.7143 -> 1258.3334
$ws.Cells.Item(102, 12).Value = 0  # L102: 1000 -> 0
$ws.Cells.Item(102, 13).Value = 363.6666  # M102: 851.2857 -> 363.6666
$ws.Cells.Item(102, 14).ClearContents()  # N102: -4244 -> (removed)
$ws.Cells.Item(132, 8).Value = 2593.25  # H132: 2631.0625 -> 2593.25
$ws.Cells.Item(132, 9).Value = 2215.8484  # I132: 2231.4243 -> 2215.8484
$ws.Cells.Item(132, 10).Value = 3423.5334  # J132: 3510.2666 -> 3423.5334
$ws.Cells.Item(132, 11).Value = 6647.5452  # K132: 6694.2729 -> 6647.5452
$ws.Cells.Item(132, 12).Value = 10270.6002  # L132: 10530.7998 -> 10270.6002
$ws.Cells.Item(132, 13).Value = -4117.5452  # M132: -4164.2729 -> -4117.5452
$ws.Cells.Item(132, 14).Value = -15330.6002  # N132: -15590.7998 -> -15330.6002
$ws.Cells.Item(136, 8).Value = 3134.5557  # H136: 3585.5715 -> 3134.5557
$ws.Cells.Item(136, 9).Value = 2078  # I136: 1999.5 -> 2078
$ws.Cells.Item(136, 10).Value = 3979.8  # J136: 4220 -> 3979.8
$ws.Cells.Item(136, 11).Value = 6234  # K136: 5998.5 -> 6234
$ws.Cells.Item(136, 12).Value = 11939.4  # L136: 12660 -> 11939.4
$ws.Cells.Item(136, 13).Value = -3684  # M136: -3448.5 -> -3684
$ws.Cells.Item(136, 14).Value = -17039.4  # N136: -17760 -> -17039.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1152.75  # H107: 1978 -> 1152.75
$ws.Cells.Item(107, 9).Value = 1112.091  # I107: 1963.3334 -> 1112.091
$ws.Cells.Item(107, 10).Value = 1600  # J107: 2000 -> 1600
$ws.Cells.Item(107, 11).Value = 1112.091  # K107: 1963.3334 -> 1112.091
$ws.Cells.Item(107, 12).Value = 1600  # L107: 2000 -> 1600
$ws.Cells.Item(107, 13).Value = 807.9090000000001  # M107: -43.33339999999998 -> 807.9090000000001
$ws.Cells.Item(107, 14).Value = -5440  # N107: -5840 -> -5440
$ws.Cells.Item(134, 8).Value = 1339.0769  # H134: 1384.7273 -> 1339.0769
$ws.Cells.Item(134, 9).Value = 1159  # I134: 1212.5714 -> 1159
$ws.Cells.Item(134, 10).Value = 3500  # J134: 5000 -> 3500
$ws.Cells.Item(134, 11).Value = 3477  # K134: 3637.7142 -> 3477
$ws.Cells.Item(134, 12).Value = 10500  # L134: 15000 -> 10500
$ws.Cells.Item(134, 13).Value = -942  # M134: -1102.7142 -> -942
$ws.Cells.Item(134, 14).Value = -15570  # N134: -20070 -> -15570

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 671.7143  # H22: 1590.6 -> 671.7143
$ws.Cells.Item(22, 9).Value = 200  # I22: 5000 -> 200
$ws.Cells.Item(22, 10).Value = 750.3333  # J22: 738.25 -> 750.3333
$ws.Cells.Item(22, 11).Value = 200  # K22: 5000 -> 200
$ws.Cells.Item(22, 12).Value = 750.3333  # L22: 738.25 -> 750.3333
$ws.Cells.Item(22, 13).Value = 150  # M22: -4650 -> 150
$ws.Cells.Item(22, 14).Value = -1450.3333  # N22: -1438.25 -> -1450.3333
$ws.Cells.Item(23, 8).Value = 34949  # H23: 36172.5 -> 34949
$ws.Cells.Item(23, 10).Value = 34949  # J23: 36172.5 -> 34949
$ws.Cells.Item(23, 12).Value = 34949  # L23: 36172.5 -> 34949
$ws.Cells.Item(23, 14).Value = -35429  # N23: -36652.5 -> -35429
$ws.Cells.Item(27, 8).Value = 34949  # H27: 36172.5 -> 34949
$ws.Cells.Item(27, 10).Value = 34949  # J27: 36172.5 -> 34949
$ws.Cells.Item(27, 12).Value = 34949  # L27: 36172.5 -> 34949
$ws.Cells.Item(27, 14).Value = -35333  # N27: -36556.5 -> -35333
$ws.Cells.Item(31, 8).Value = 8125.75  # H31: 8092.591 -> 8125.75
$ws.Cells.Item(31, 9).Value = 3273.4814  # I31: 3209.074 -> 3273.4814
$ws.Cells.Item(31, 10).Value = 15832.294  # J31: 15848.765 -> 15832.294
$ws.Cells.Item(31, 11).Value = 3273.4814  # K31: 3209.074 -> 3273.4814
$ws.Cells.Item(31, 12).Value = 15832.294  # L31: 15848.765 -> 15832.294
$ws.Cells.Item(31, 13).Value = -2978.4814  # M31: -2914.074 -> -2978.4814
$ws.Cells.Item(31, 14).Value = -16422.294  # N31: -16438.765 -> -16422.294
$ws.Cells.Item(34, 8).Value = 8125.75  # H34: 8092.591 -> 8125.75
$ws.Cells.Item(34, 9).Value = 3273.4814  # I34: 3209.074 -> 3273.4814
$ws.Cells.Item(34, 10).Value = 15832.294  # J34: 15848.765 -> 15832.294
$ws.Cells.Item(34, 11).Value = 3273.4814  # K34: 3209.074 -> 3273.4814
$ws.Cells.Item(34, 12).Value = 15832.294  # L34: 15848.765 -> 15832.294
$ws.Cells.Item(34, 13).Value = -3071.4814  # M34: -3007.074 -> -3071.4814
$ws.Cells.Item(34, 14).Value = -16236.294  # N34: -16252.765 -> -16236.294

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(21, 8).Value = 750  # H21: 100 -> 750
$ws.Cells.Item(21, 9).Value = 583.3333  # I21: 100 -> 583.3333
$ws.Cells.Item(21, 10).Value = 1500  # J21: 0 -> 1500
$ws.Cells.Item(21, 11).Value = 1749.9999  # K21: 300 -> 1749.9999
$ws.Cells.Item(21, 12).Value = 4500  # L21: 0 -> 4500
$ws.Cells.Item(21, 13).Value = -1576.9999  # M21: -127 -> -1576.9999
$ws.Cells.Item(21, 14).Value = -4846  # N21: None -> -4846
$ws.Cells.Item(98, 8).Value = 857.6429000000001  # H98: 435878.75 -> 857.6429000000001
$ws.Cells.Item(98, 9).Value = 694.125  # I98: 1200.4 -> 694.125
$ws.Cells.Item(98, 10).Value = 1075.6666  # J98: 770246.7 -> 1075.6666
$ws.Cells.Item(98, 11).Value = 2082.375  # K98: 3601.2 -> 2082.375
$ws.Cells.Item(98, 12).Value = 3226.9998  # L98: 2310740.1 -> 3226.9998
$ws.Cells.Item(98, 13).Value = -584.375  # M98: -2103.2 -> -584.375
$ws.Cells.Item(98, 14).Value = -6222.9998  # N98: -2313736.1 -> -6222.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2344.25  # H40: 2675 -> 2344.25
$ws.Cells.Item(40, 9).Value = 2050.6667  # I40: 2400 -> 2050.6667
$ws.Cells.Item(40, 11).Value = 2050.6667  # K40: 2400 -> 2050.6667
$ws.Cells.Item(40, 13).Value = -1914.6667  # M40: -2264 -> -1914.6667
$ws.Cells.Item(46, 8).Value = 391339.5  # H46: 782299.8 -> 391339.5
$ws.Cells.Item(46, 9).Value = 300  # I46: 0 -> 300
$ws.Cells.Item(46, 10).Value = 489099.38  # J46: 782299.8 -> 489099.38
$ws.Cells.Item(46, 11).Value = 300  # K46: 0 -> 300
$ws.Cells.Item(46, 12).Value = 489099.38  # L46: 782299.8 -> 489099.38
$ws.Cells.Item(46, 13).Value = -112  # M46: None -> -112
$ws.Cells.Item(46, 14).Value = -489475.38  # N46: -782675.8 -> -489475.38
